# Replacing test data files with latest.
# The "decision notice date" row (row 2) is updated from 07/2019 to 01/2020,
# and the active selection moves to the cell that was just edited (D2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Month: 07 -> 01
$ws.Range("C2").Value = "01"

# Year: 2019 -> 2020 (stored as text, matching the column's text format)
$ws.Range("D2").Value = "2020"

# Leave the selection on the cell that was last edited.
$ws.Range("D2").Select() | Out-Null
